# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" stats sheet with a newer data pull.
# A handful of neighbouring countries swapped rank (so the label shown
# in a given row changes) and their metrics were refreshed; a few other
# rows only got refreshed metrics with the same country label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 - Arabia Saudita (same country, refreshed numbers)
$ws.Range("B25").Value = 17522
$ws.Range("C25").Value = 1223
$ws.Range("D25").Value = 2357
$ws.Range("E25").Value = 15026
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 139

# Row 41 - Dinamarca (same country, refreshed numbers)
$ws.Range("B41").Value = 8575
$ws.Range("C41").Value = 130
$ws.Range("D41").Value = 5805
$ws.Range("E41").Value = 2348
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 422

# Row 58 - Moldavia (same country, refreshed numbers)
$ws.Range("D58").Value = 895
$ws.Range("E58").Value = 2315

# Row 62 - Kazajistan (same country, refreshed numbers)
$ws.Range("E62").Value = 1991
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 25

# Row 68 - Uzbekistan (same country, refreshed numbers)
$ws.Range("B68").Value = 1869
$ws.Range("C68").Value = 7
$ws.Range("D68").Value = 765
$ws.Range("E68").Value = 1096

# Rows 74-76 - Ghana/Afganistan/Camerun re-ranked: Camerun moves up to 74,
# Ghana to 75, Afganistan to 76, each with refreshed numbers.
$ws.Range("A74").Value = "Camerun"
$ws.Range("B74").Value = 1621
$ws.Range("C74").Value = 103
$ws.Range("D74").Value = 786
$ws.Range("E74").Value = 779
$ws.Range("F74").Value = 12
$ws.Range("G74").Value = 3
$ws.Range("H74").Value = 56

$ws.Range("A75").Value = "Ghana"
$ws.Range("B75").Value = 1550
$ws.Range("C75").Value = 271
$ws.Range("D75").Value = 155
$ws.Range("E75").Value = 1384
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 11

$ws.Range("A76").Value = "Afganistan"
$ws.Range("B76").Value = 1531
$ws.Range("C76").Value = 68
$ws.Range("D76").Value = 207
$ws.Range("E76").Value = 1274
$ws.Range("F76").Value = 7
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = 50

# Rows 101-103 - Burkina Faso/Honduras/Senegal re-ranked: Senegal moves up
# to 101, Burkina Faso to 102, Honduras to 103, each with refreshed numbers.
$ws.Range("A101").Value = "Senegal"
$ws.Range("B101").Value = 671
$ws.Range("C101").Value = 57
$ws.Range("D101").Value = 283
$ws.Range("E101").Value = 380
$ws.Range("F101").Value = 1
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 8

$ws.Range("A102").Value = "Burkina Faso"
$ws.Range("B102").Value = 629
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 442
$ws.Range("E102").Value = 146
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 41

$ws.Range("A103").Value = "Honduras"
$ws.Range("B103").Value = 627
$ws.Range("C103").Value = 36
$ws.Range("D103").Value = 65
$ws.Range("E103").Value = 503
$ws.Range("F103").Value = 10
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 59

# Rows 107-108 - Guatemala/Sri Lanka re-ranked: Sri Lanka moves up to 107,
# Guatemala to 108, each with refreshed numbers.
$ws.Range("A107").Value = "Sri Lanka"
$ws.Range("B107").Value = 477
$ws.Range("C107").Value = 25
$ws.Range("D107").Value = 120
$ws.Range("E107").Value = 350
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 7

$ws.Range("A108").Value = "Guatemala"
$ws.Range("B108").Value = 473
$ws.Range("C108").Value = 43
$ws.Range("D108").Value = 45
$ws.Range("E108").Value = 415
$ws.Range("F108").Value = 5
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 13

# Row 167 - Nepal (same country, refreshed numbers)
$ws.Range("B167").Value = 52
$ws.Range("C167").Value = 3
$ws.Range("E167").Value = 36
